$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the new time-log entry for row 43 (date 2014-09-15)
# (Interruption is written first so the shared formula in E43 picks it up
# correctly when it recalculates.)
$ws.Range("D43").Value = 15
$ws.Range("A43").Value = 41897
$ws.Range("B43").Value = 0.8847222222222223
$ws.Range("C43").Value = 0.9819444444444444
$ws.Range("F43").Value = "Testing"

$excel.Calculate()

# Update the selected cell to match the author's final cursor position
$ws.Range("A44").Select()
